# Luan Phuc Duc - add 79 new "hoi chieu tai Phuc Duc" combination rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    "Sát Phá Tham hội chiếu tại Phúc Đức",
    "Tử Phủ Vũ Tướng hội chiếu tại Phúc Đức",
    "Cơ Nguyệt Đồng Lương hội chiếu tại Phúc Đức",
    "Cự Nhật hội chiếu tại Phúc Đức",
    "Kình Đà hội chiếu tại Phúc Đức",
    "Xương Khúc hội chiếu tại Phúc Đức",
    "Hoả Linh hội chiếu tại Phúc Đức",
    "Không Kiếp hội chiếu tại Phúc Đức",
    "Quang Quý hội chiếu tại Phúc Đức",
    "Tả Hữu hội chiếu tại Phúc Đức",
    "Song Hao hội chiếu tại Phúc Đức",
    "Tang Hổ hội chiếu tại Phúc Đức",
    "Khốc Hư hội chiếu tại Phúc Đức",
    "Hình Riêu hội chiếu tại Phúc Đức",
    "Thai Toạ hội chiếu tại Phúc Đức",
    "Đào Hồng hội chiếu tại Phúc Đức",
    "Ấn Phù hội chiếu tại Phúc Đức",
    "Song Hao Quyền Lộc Kiếp Hoả hội chiếu tại Phúc Đức",
    "Tử Phủ Vũ Tướng Xương Khúc Khôi Việt Tả Hữu Khoa Quyền Lộc Long hội chiếu tại Phúc Đức",
    "Tử Phủ Vũ Tướng Tả Hữu Khoa Quyền Lộc Long Phượng hội chiếu tại Phúc Đức",
    "Tử Khúc Phá Dương Đà hội chiếu tại Phúc Đức",
    "Cơ Nguyệt Đồng Lương gặp Xương Khúc Tả Hữu hội chiếu tại Phúc Đức",
    "Cơ Nguyệt Đồng Lương Khoa Tả Hữu Quang Quý Quan Phúc hội chiếu tại Phúc Đức",
    "Sát Quyền hội chiếu tại Phúc Đức",
    "Lộc Mã hội chiếu tại Phúc Đức",
    "Kiếp Hư Hao Quyền hội chiếu tại Phúc Đức",
    "Tuế Hổ Phù Xương Khúc hội chiếu tại Phúc Đức",
    "Xương Khúc Tấu Long Phượng hội chiếu tại Phúc Đức",
    "Đào Hồng Riêu Tấu Cơ Vũ hội chiếu tại Phúc Đức",
    "Binh Hình Tướng Ấn hội chiếu tại Phúc Đức",
    "Hổ Tấu hội chiếu tại Phúc Đức",
    "Hình Riêu Y hội chiếu tại Phúc Đức",
    "Mã Hỏa Linh hội chiếu tại Phúc Đức",
    "Thai Tọa Hồng Đào hội chiếu tại Phúc Đức",
    "Tả Hữu Không Kiếp hội chiếu tại Phúc Đức",
    "Tả Hữu Binh Tướng hội chiếu tại Phúc Đức",
    "Đào Quyền hội chiếu tại Phúc Đức",
    "Đào Hồng Tả Cái Hữu hội chiếu tại Phúc Đức",
    "Quan Phúc Quang Tấu hội chiếu tại Phúc Đức",
    "Đào Hồng Xương Khúc hội chiếu tại Phúc Đức",
    "Đào Tử Phủ hội chiếu tại Phúc Đức",
    "Hổ Kình Sát hội chiếu tại Phúc Đức",
    "Hổ Tang Kiếp hội chiếu tại Phúc Đức",
    "Hổ Tang Không Kiếp hội chiếu tại Phúc Đức",
    "Âm Dương Lương hội chiếu tại Phúc Đức",
    "Cơ Lương Gia Hội hội chiếu tại Phúc Đức",
    "Nhật Chiếu Lôi Môn hội chiếu tại Phúc Đức",
    "Tả Hữu Xương Khúc hội chiếu tại Phúc Đức",
    "Tham Linh Triều Viên hội chiếu tại Phúc Đức",
    "Cự Hỏa Linh hội chiếu tại Phúc Đức",
    "Hồng Đào Kỵ hội chiếu tại Phúc Đức",
    "Cự Đồng Hình hội chiếu tại Phúc Đức",
    "Lương Phá hội chiếu tại Phúc Đức",
    "Khôi Việt hội chiếu tại Phúc Đức",
    "Tham Vũ Hỏa hội chiếu tại Phúc Đức",
    "Cự Kỵ hội chiếu tại Phúc Đức",
    "Kình Đà Hỏa hội chiếu tại Phúc Đức",
    "Kình Đà Không Kiếp hội chiếu tại Phúc Đức",
    "Đào Hồng Hỷ hội chiếu tại Phúc Đức",
    "Đào Hồng Kiếp Sát hội chiếu tại Phúc Đức",
    "Đào Hồng Kỵ hội chiếu tại Phúc Đức",
    "Đào Hồng Tả Phù Hữu Bật hội chiếu tại Phúc Đức",
    "Đào Hồng Tả Phù Hữu Bật Khoa Quyền Lộc hội chiếu tại Phúc Đức",
    "Đào Hồng Tả Phù Hữu Bật Khoa Quyền Lộc Long Trì Phượng Các hội chiếu tại Phúc Đức",
    "Lưu Hà Kiếp Sát hội chiếu tại Phúc Đức",
    "Phục Không Kiếp hội chiếu tại Phúc Đức",
    "Xương Khúc Khôi Việt hội chiếu tại Phúc Đức",
    "Khoa Quyền hội chiếu tại Phúc Đức",
    "Lộc Quyền hội chiếu tại Phúc Đức",
    "Tướng Binh Đào Hồng hội chiếu tại Phúc Đức",
    "Lương Khốc Tuế hội chiếu tại Phúc Đức",
    "Tướng Binh hội chiếu tại Phúc Đức",
    "Cự Tang hội chiếu tại Phúc Đức",
    "Cự Hỏa hội chiếu tại Phúc Đức",
    "Phá Hình Kỵ hội chiếu tại Phúc Đức",
    "Tang Trực Tuế hội chiếu tại Phúc Đức",
    "Tuế Xương Khúc hội chiếu tại Phúc Đức",
    "Tử Vi Tả Hữu hội chiếu tại Phúc Đức",
    "Tử Vi Tang Tả Hữu hội chiếu tại Phúc Đức"
)

$startRow = 4244
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newValues[$i]
    $ws.Cells.Item($r, 2).Value = $newValues[$i]
}

$lastRow = $startRow + $newValues.Length - 1

# New conditional formatting rule ("Highlight duplicate values") for the freshly added A column rows,
# matching the existing red-on-red duplicate style used throughout the sheet, placed at top priority.
$newRange = $ws.Range("A" + $startRow + ":A" + $lastRow)
$fc = $newRange.FormatConditions.AddUniqueValues()
$fc.DupeUnique = 1
$fc.Font.Color = 393372
$fc.Interior.Color = 13551615
$fc.SetFirstPriority()

# Match the author's final selection/scroll state: viewport scrolled near the bottom, with the
# newly added column B cells selected.
[void]$excel.Goto($ws.Range("A4212"), $true)
[void]$ws.Range("B" + $startRow + ":B" + $lastRow).Select()

Write-Host ("Added " + $newValues.Length + " new rows (" + $startRow + " to " + $lastRow + ") and conditional formatting.")
